# Apply stimulus renaming (face -> book) and expand shorthand answer codes
# (y/b/r -> left/center/right) across the trial-sequence worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 361

# 1) Column L ("correct_ans"): expand single-letter codes to full words.
$answerMap = @{ "y" = "left"; "b" = "center"; "r" = "right" }
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 12)
    $cur = $cell.Value()
    if ($answerMap.ContainsKey($cur)) {
        $cell.Value = $answerMap[$cur]
    }
}

# 2) Columns A-D (prompt/correct/distractor image files): rename the
#    "face" stimulus category to "book" (e.g. face//face_40.jpg -> book//book_40.jpg).
$imageCols = 1,2,3,4
foreach ($c in $imageCols) {
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $c)
        $cur = $cell.Value()
        if ($cur -like "face//*") {
            $cell.Value = $cur -replace "face", "book"
        }
    }
}
